$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (bold)
$ws.Range("A1").Value = "Akronym"
$ws.Range("B1").Value = "Beskrivelse"
$ws.Range("A1:B1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "CD-ROM"
$ws.Range("B2").Value = "Compact disk read only memory"
$ws.Range("A3").Value = "RAM"

# Column widths to match "bestFit" autofit sizing
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Select the last-edited cell, matching the saved selection in the sheet view
$ws.Range("B3").Select() | Out-Null
